# Updates "想去人数" (col F) and "最低票价" (col G) figures across the four
# worksheets ("展览", "演出", "本地生活", "全部类型") to match the refreshed
# crawl data, per the commit's regenerated output.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 1179
$ws.Range("G3").Value = 40.5
$ws.Range("F4").Value = 1253
$ws.Range("F6").Value = 173
$ws.Range("F7").Value = 540
$ws.Range("G7").Value = 54
$ws.Range("F8").Value = 327
$ws.Range("F10").Value = 1260
$ws.Range("F11").Value = 28652
$ws.Range("F12").Value = 3491
$ws.Range("F13").Value = 36
$ws.Range("F14").Value = 251
$ws.Range("F15").Value = 474
$ws.Range("F16").Value = 26
$ws.Range("G16").Value = 45
$ws.Range("F18").Value = 11
$ws.Range("F20").Value = 616
$ws.Range("F21").Value = 273
$ws.Range("F22").Value = 266
$ws.Range("F23").Value = 351
$ws.Range("F25").Value = 53
$ws.Range("F26").Value = 660
$ws.Range("G26").Value = 54
$ws.Range("F27").Value = 207
$ws.Range("F28").Value = 99
$ws.Range("F29").Value = 536
$ws.Range("F30").Value = 71
$ws.Range("F31").Value = 34
$ws.Range("F32").Value = 620
$ws.Range("F35").Value = 1

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F2").Value = 3
$ws.Range("F4").Value = 18
$ws.Range("F6").Value = 381
$ws.Range("F7").Value = 840
$ws.Range("F9").Value = 88
$ws.Range("F10").Value = 269
$ws.Range("F22").Value = 4242

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 255
$ws.Range("F4").Value = 1174

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 255
$ws.Range("F4").Value = 1174
$ws.Range("F5").Value = 3
$ws.Range("F6").Value = 18
$ws.Range("F7").Value = 381
$ws.Range("F9").Value = 840
$ws.Range("F10").Value = 1179
$ws.Range("G10").Value = 40.5
$ws.Range("F11").Value = 1253
$ws.Range("F12").Value = 173
$ws.Range("F13").Value = 540
$ws.Range("G13").Value = 54
$ws.Range("F14").Value = 327
$ws.Range("F16").Value = 52
$ws.Range("F17").Value = 1260
$ws.Range("F18").Value = 88
$ws.Range("F19").Value = 88
$ws.Range("F20").Value = 269
$ws.Range("F25").Value = 251
$ws.Range("F28").Value = 474
$ws.Range("F29").Value = 26
$ws.Range("G29").Value = 45
$ws.Range("F33").Value = 616
$ws.Range("F34").Value = 273
$ws.Range("F35").Value = 351
$ws.Range("F37").Value = 53
$ws.Range("F38").Value = 660
$ws.Range("G38").Value = 54
$ws.Range("F40").Value = 207
$ws.Range("F41").Value = 99
$ws.Range("F44").Value = 71
$ws.Range("F45").Value = 34
$ws.Range("F46").Value = 620
$ws.Range("F50").Value = 1
